# Updated mean summary results
#
# The source data table (common_name_e, day_night, habitat_type, accel, sem)
# contained one extra summary row per (species, day_night) group where the
# habitat_type (column C) was blank/NA. Those 8 rows are removed, which
# shifts all subsequent rows up.
#
# Rows to delete (1-based, as found in the original sheet) are the ones
# whose column C is empty: 7, 13, 19, 25, 31, 37, 43, 49.
# Delete from the bottom up so earlier row numbers remain valid targets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(49, 43, 37, 31, 25, 19, 13, 7)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
